$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.121.69"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "2.370.04"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'304.10"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'95.51"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.503"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("D10").Value = "'34.37"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").Value = "'18.60"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").Value = "2.736.58"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").Value = "2.361.77"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "43.109.08"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "'12.00"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").Value = "'68.16"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'235.49"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("E24").Value = "  -2.16%  "

$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'24.60"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("D29").Value = "'9.37"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("D30").Value = "'32.49"
$ws.Range("E30").Value = "  +2.90%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").Value = "'5.03"
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("D33").Value = "'17.55"
$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("D34").Value = "'0.0729"
$ws.Range("E34").Value = "  +3.93%  "

$ws.Range("D35").Value = "'0.107"
$ws.Range("E35").Value = "  +6.33%  "

$ws.Range("E36").Value = "  +1.72%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'126.96"
$ws.Range("E37").Value = "  -9.88%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.35"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.85"
$ws.Range("E39").Value = "  +3.66%  "

$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "'20.79"
$ws.Range("E42").Value = "  -7.36%  "

$ws.Range("D43").Value = "1.936.63"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("E45").Value = "  +3.99%  "

$ws.Range("E46").Value = "  -7.79%  "

$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").Value = "2.596.10"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("D49").Value = "'1.52"
$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("D50").Value = "'71.91"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("D51").Value = "'1.13"
$ws.Range("E51").Value = "  +0.97%  "
